# The deck ships two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (default Office palette) -- only
#                             wired to the Notes Master, not reachable via the
#                             PowerPoint object model.
#   ppt/theme/theme2.xml  -> "Integral" theme -- the theme actually driving
#                             the Slide Master / presentation design, and the
#                             only theme object exposed through COM
#                             (SlideMaster.Theme / Design.SlideMaster.Theme /
#                             NotesMaster.Theme / HandoutMaster.Theme all
#                             resolve to this same live theme).
#
# The authored edit swaps the two themes' content: the deck's working theme
# becomes the plain "Office" color palette (dk2/lt2/accent1-6/hlink/folHlink)
# while the font scheme and format scheme (already identical, Arial-based
# "Office" scheme, between the two parts) stay the same. Apply that palette
# swap to the live theme via its ThemeColorScheme, color by color (the
# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink). RGB values use the standard VBA RGB(r,g,b) = r + g*256 + b*65536
# encoding.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

$cs.Item(1).RGB  = 0          # dk1      -> 000000
$cs.Item(2).RGB  = 16777215   # lt1      -> FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      -> 44546A
$cs.Item(4).RGB  = 15132391   # lt2      -> E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  -> ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  -> A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  -> FFC000
$cs.Item(9).RGB  = 12874308   # accent5  -> 4472C4
$cs.Item(10).RGB = 4697456    # accent6  -> 70AD47
$cs.Item(11).RGB = 12673797   # hlink    -> 0563C1
$cs.Item(12).RGB = 7491477    # folHlink -> 954F72
